$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44441, 0, 12, 110.5481345002303),
    @(44442, 0, 7, 64.486411791801),
    @(44443, 1, 6, 55.27406725011516),
    @(44444, 2, 7, 64.486411791801),
    @(44445, 0, 7, 64.486411791801),
    @(44446, 2, 5, 46.06172270842929),
    @(44447, 1, 6, 55.27406725011516),
    @(44448, 1, 7, 64.486411791801)
)

$startRow = 367

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    # Match the formatting used in column A for the preceding rows (date-style, bordered, bold, centered)
    $ws.Range("A366").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
